$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update Maximo (C2) with new minimized average time result ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 641.831360147043

# --- Sheet "Solucion": randomized re-assignment of Pedido -> Salida pairing ---
$wsSolucion = $wb.Worksheets.Item("Solucion")

$pedidos = @(
    "Pedido_16",
    "Pedido_28",
    "Pedido_35",
    "Pedido_5",
    "Pedido_7",
    "Pedido_15",
    "Pedido_22",
    "Pedido_33",
    "Pedido_6",
    "Pedido_38",
    "Pedido_2",
    "Pedido_3",
    "Pedido_17",
    "Pedido_40",
    "Pedido_25",
    "Pedido_14",
    "Pedido_18",
    "Pedido_37",
    "Pedido_19",
    "Pedido_4",
    "Pedido_13",
    "Pedido_9",
    "Pedido_23",
    "Pedido_11",
    "Pedido_12",
    "Pedido_39",
    "Pedido_10",
    "Pedido_27",
    "Pedido_36",
    "Pedido_32",
    "Pedido_1",
    "Pedido_26",
    "Pedido_20",
    "Pedido_30",
    "Pedido_24",
    "Pedido_31",
    "Pedido_21",
    "Pedido_29",
    "Pedido_34",
    "Pedido_8"
)

$salidas = @(
    "S001",
    "S025",
    "S029",
    "S005",
    "S026",
    "S002",
    "S030",
    "S006",
    "S027",
    "S003",
    "S007",
    "S031",
    "S028",
    "S004",
    "S008",
    "S032",
    "S033",
    "S009",
    "S013",
    "S037",
    "S010",
    "S034",
    "S038",
    "S014",
    "S035",
    "S011",
    "S039",
    "S015",
    "S036",
    "S012",
    "S016",
    "S040",
    "S017",
    "S021",
    "S018",
    "S022",
    "S019",
    "S023",
    "S020",
    "S024"
)

for ($i = 0; $i -lt $pedidos.Length; $i++) {
    $row = 2 + $i
    $wsSolucion.Cells.Item($row, 1).Value = $pedidos[$i]
    $wsSolucion.Cells.Item($row, 2).Value = $salidas[$i]
}

# --- Sheet "Metricas": update Tiempo per Zona with new minimized average time results ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 641.831360147043
$wsMetricas.Range("B3").Value = 450.790896313115

